$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "as of" timestamp in the title row (A1) of every worksheet.
#    13:35 -> 17:39 (same date, 30/04/2021).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MD410 Attendance")
$ws2 = $wb.Worksheets.Item("410E Attendance")
$ws3 = $wb.Worksheets.Item("410W Attendance")
$ws4 = $wb.Worksheets.Item("410E Voting")
$ws5 = $wb.Worksheets.Item("410W Voting")

$ws1.Cells.Item(1,1).Value = "MD410 Registrees as of 30/04/2021 17:39"
$ws2.Cells.Item(1,1).Value = "410E Registrees as of 30/04/2021 17:39"
$ws3.Cells.Item(1,1).Value = "410W Registrees as of 30/04/2021 17:39"
$ws4.Cells.Item(1,1).Value = "410E Voting details as of 30/04/2021 17:39"
$ws5.Cells.Item(1,1).Value = "410W Voting details as of 30/04/2021 17:39"

# ---------------------------------------------------------------------------
# 2. Insert a new registree row on "MD410 Attendance" at row 152
#    (alphabetically between "Pantoleon" and "Pillay"), pushing every
#    following row down by one.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(152).Insert()

# Carry over the row-level formatting (height) and cell formatting (style,
# borders) from the row immediately above, then overwrite with the new
# registree's values.
$ws1.Range("A151:F151").Copy()
$ws1.Range("A152:F152").PasteSpecial(-4122)
$ws1.Rows.Item(152).RowHeight = 25

$ws1.Cells.Item(152,1).Value = "Piater"
$ws1.Cells.Item(152,2).Value = "Ivan"
$ws1.Cells.Item(152,3).Value = "East London Port Rex"
$ws1.Cells.Item(152,4).Value = "No"
$ws1.Cells.Item(152,5).Value = "No"
$ws1.Cells.Item(152,6).Value = "410E"

# ---------------------------------------------------------------------------
# 3. Bump the "Number of attendees" summary row to reflect the new total
#    (246 -> 247). This row was pushed from 249 to 250 by the insert above.
#    "Number of voters: 98" (now row 251) is unaffected.
# ---------------------------------------------------------------------------
$ws1.Cells.Item(250,1).Value = "Number of attendees: 247"
